$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 270.0909
$ws.Range("I33").Value = 254.2963
$ws.Range("K33").Value = 254.2963
$ws.Range("M33").Value = -25.2963
$ws.Range("H43").Value = 2004.0588
$ws.Range("I43").Value = 844.5
$ws.Range("J43").Value = 2158.6667
$ws.Range("K43").Value = 844.5
$ws.Range("L43").Value = 2158.6667
$ws.Range("M43").Value = -775.5
$ws.Range("N43").Value = -2296.6667
$ws.Range("H129").Value = 869.7143
$ws.Range("J129").Value = 1030.4706
$ws.Range("L129").Value = 3091.4118
$ws.Range("N129").Value = -13091.4118
$ws.Range("H132").Value = 540896
$ws.Range("I132").Value = 2258.883
$ws.Range("J132").Value = 3503400.2
$ws.Range("K132").Value = 6776.648999999999
$ws.Range("L132").Value = 10510200.6
$ws.Range("M132").Value = -4246.648999999999
$ws.Range("N132").Value = -10515260.6
$ws.Range("H133").Value = 29780
$ws.Range("J133").Value = 29780
$ws.Range("L133").Value = 29780
$ws.Range("N133").Value = -39900
$ws.Range("H137").Value = 4767046.5
$ws.Range("I137").Value = 6254249
$ws.Range("J137").Value = 7998
$ws.Range("K137").Value = 18762747
$ws.Range("L137").Value = 23994
$ws.Range("M137").Value = -18760197
$ws.Range("N137").Value = -29094
$ws.Range("H138").Value = 3511505.8
$ws.Range("I138").Value = 2168.1143
$ws.Range("J138").Value = 9094543
$ws.Range("K138").Value = 6504.342900000001
$ws.Range("L138").Value = 27283629
$ws.Range("M138").Value = -1364.342900000001
$ws.Range("N138").Value = -27293909

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14664.47
$ws.Range("I32").Value = 10356.029
$ws.Range("J32").Value = 24254.227
$ws.Range("K32").Value = 10356.029
$ws.Range("L32").Value = 24254.227
$ws.Range("M32").Value = -10069.029
$ws.Range("N32").Value = -24828.227
$ws.Range("H61").Value = 27835738
$ws.Range("I61").Value = 34518750
$ws.Range("J61").Value = 148957.14
$ws.Range("K61").Value = 34518750
$ws.Range("L61").Value = 148957.14
$ws.Range("M61").Value = -34518538
$ws.Range("N61").Value = -149381.14
$ws.Range("H132").Value = 13215058
$ws.Range("I132").Value = 18557592
$ws.Range("J132").Value = 101566.18
$ws.Range("K132").Value = 55672776
$ws.Range("L132").Value = 304698.54
$ws.Range("M132").Value = -55670246
$ws.Range("N132").Value = -309758.54
$ws.Range("H136").Value = 27835738
$ws.Range("I136").Value = 34518750
$ws.Range("J136").Value = 148957.14
$ws.Range("K136").Value = 103556250
$ws.Range("L136").Value = 446871.42
$ws.Range("M136").Value = -103553700
$ws.Range("N136").Value = -451971.42

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7248850
$ws.Range("I134").Value = 2189.0938
$ws.Range("J134").Value = 23812646
$ws.Range("K134").Value = 6567.2814
$ws.Range("L134").Value = 71437938
$ws.Range("M134").Value = -4032.2814
$ws.Range("N134").Value = -71443008

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 24426.092
$ws.Range("I132").Value = 1433.1111
$ws.Range("J132").Value = 60944.35
$ws.Range("K132").Value = 4299.3333
$ws.Range("L132").Value = 182833.05
$ws.Range("M132").Value = -1769.3333
$ws.Range("N132").Value = -187893.05
$ws.Range("H134").Value = 45763.668
$ws.Range("I134").Value = 770.1875
$ws.Range("J134").Value = 135750.62
$ws.Range("K134").Value = 2310.5625
$ws.Range("L134").Value = 407251.86
$ws.Range("M134").Value = 224.4375
$ws.Range("N134").Value = -412321.86

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2500.2
$ws.Range("J75").Value = 2800.25
$ws.Range("L75").Value = 8400.75
$ws.Range("N75").Value = -10396.75
$ws.Range("H78").Value = 2500.2
$ws.Range("J78").Value = 2800.25
$ws.Range("L78").Value = 25202.25
$ws.Range("N78").Value = -35186.25
$ws.Range("H131").Value = 968
$ws.Range("I131").Value = 770
$ws.Range("J131").Value = 1001
$ws.Range("K131").Value = 2310
$ws.Range("L131").Value = 3003
$ws.Range("M131").Value = 2730
$ws.Range("N131").Value = -13083

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1862.0333
$ws.Range("I102").Value = 1744.9131
$ws.Range("J102").Value = 2246.8572
$ws.Range("K102").Value = 1744.9131
$ws.Range("L102").Value = 2246.8572
$ws.Range("M102").Value = -122.9131
$ws.Range("N102").Value = -5490.8572
$ws.Range("H122").Value = 3131.1538
$ws.Range("I122").Value = 2685.2856
$ws.Range("J122").Value = 3651.3333
$ws.Range("K122").Value = 8055.8568
$ws.Range("L122").Value = 10953.9999
$ws.Range("M122").Value = -5605.8568
$ws.Range("N122").Value = -15853.9999
$ws.Range("H132").Value = 60103.516
$ws.Range("I132").Value = 49983.43
$ws.Range("J132").Value = 75283.64
$ws.Range("K132").Value = 149950.29
$ws.Range("L132").Value = 225850.92
$ws.Range("M132").Value = -147420.29
$ws.Range("N132").Value = -230910.92
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1806.2759
$ws.Range("I61").Value = 1863.0588
$ws.Range("J61").Value = 1725.8334
$ws.Range("K61").Value = 1863.0588
$ws.Range("L61").Value = 1725.8334
$ws.Range("M61").Value = -1661.0588
$ws.Range("N61").Value = -2129.8334
$ws.Range("H113").Value = 1806.2759
$ws.Range("I113").Value = 1863.0588
$ws.Range("J113").Value = 1725.8334
$ws.Range("K113").Value = 1863.0588
$ws.Range("L113").Value = 1725.8334
$ws.Range("M113").Value = 306.9412
$ws.Range("N113").Value = -6065.8334
$ws.Range("H122").Value = 3987.75
$ws.Range("I122").Value = 3989.1765
$ws.Range("J122").Value = 3984.2856
$ws.Range("K122").Value = 11967.5295
$ws.Range("L122").Value = 11952.8568
$ws.Range("M122").Value = -9517.529500000001
$ws.Range("N122").Value = -16852.8568
$ws.Range("H132").Value = 23856
$ws.Range("I132").Value = 1918.8334
$ws.Range("J132").Value = 102829.8
$ws.Range("K132").Value = 5756.5002
$ws.Range("L132").Value = 308489.4
$ws.Range("M132").Value = -3226.5002
$ws.Range("N132").Value = -313549.4
$ws.Range("H136").Value = 30798.584
$ws.Range("I136").Value = 22693.883
$ws.Range("J136").Value = 50481.43
$ws.Range("K136").Value = 68081.649
$ws.Range("L136").Value = 151444.29
$ws.Range("M136").Value = -65531.649
$ws.Range("N136").Value = -156544.29

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 851.5454999999999
$ws.Range("I113").Value = 914.6875
$ws.Range("J113").Value = 683.1667
$ws.Range("K113").Value = 2744.0625
$ws.Range("L113").Value = 2049.5001
$ws.Range("M113").Value = -574.0625
$ws.Range("N113").Value = -6389.5001
$ws.Range("H122").Value = 2515
$ws.Range("I122").Value = 1992.5
$ws.Range("J122").Value = 2863.3333
$ws.Range("K122").Value = 5977.5
$ws.Range("L122").Value = 8589.999899999999
$ws.Range("M122").Value = -3527.5
$ws.Range("N122").Value = -13489.9999
$ws.Range("H132").Value = 35779.793
$ws.Range("I132").Value = 28673.75
$ws.Range("J132").Value = 47407.863
$ws.Range("K132").Value = 86021.25
$ws.Range("L132").Value = 142223.589
$ws.Range("M132").Value = -83491.25
$ws.Range("N132").Value = -147283.589
